$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.816.10"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.764.38"
$ws.Range("E3").Value = "  -3.89%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.21"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.46"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.761.19"
$ws.Range("E7").Value = "  -3.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.636"
$ws.Range("E8").Value = "  -5.72%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.724"
$ws.Range("E10").Value = "  -4.72%  "
$ws.Range("E11").Value = "  -8.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.18"
$ws.Range("E12").Value = "  +6.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  -8.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.83"
$ws.Range("E14").Value = "  -7.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.366.26"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.768.08"
$ws.Range("E16").Value = "  -3.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.70"
$ws.Range("E17").Value = "  -5.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.07"
$ws.Range("E18").Value = "  -7.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.127"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("E20").Value = "  -7.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.649.24"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "418.20"
$ws.Range("E22").Value = "  -5.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.73"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.88"
$ws.Range("E24").Value = "  -5.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.08"
$ws.Range("E25").Value = "  -6.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.12"
$ws.Range("E26").Value = "  -5.46%  "
$ws.Range("E27").Value = "  -8.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.85"
$ws.Range("E28").Value = "  -5.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.09"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.74"
$ws.Range("E30").Value = "  -7.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.48"
$ws.Range("E31").Value = "  -5.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.37"
$ws.Range("E32").Value = "  -15.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.69"
$ws.Range("E33").Value = "  -8.16%  "
$ws.Range("E34").Value = "  -5.31%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.10"
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "625.73"
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "44.92"
$ws.Range("E37").Value = "  -6.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0893"
$ws.Range("E38").Value = "  -12.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.408"
$ws.Range("E39").Value = "  -6.54%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.142"
$ws.Range("E42").Value = "  -3.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.11"
$ws.Range("E43").Value = "  -7.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0447"
$ws.Range("E44").Value = "  -6.03%  "
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("E46").Value = "  -11.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.36"
$ws.Range("E47").Value = "  -9.68%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.843.42"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.137"
$ws.Range("E49").Value = "  -5.80%  "
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.13"
$ws.Range("E51").Value = "  -4.86%  "
